$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell Y1 = "22" (week 22), stored as text like the other week-number
# headers (D1..X1 are inlineStr "1".."21"). Build it as a TEXT() formula in a
# scratch cell, then paste-special just the value into Y1 so it lands as a
# shared string using the same style as X1 (s="1"), without creating a new
# number-format style.
$ws.Range("Z1").Formula = "=TEXT(22,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("Y1").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Week 22 counts for each facility row (only rows that already reported data
# for week 21 / other weeks get a week-22 cell -- matches the sparse source
# data).
$ws.Range("Y2").Value = 0
$ws.Range("Y5").Value = 0
$ws.Range("Y6").Value = 27
$ws.Range("Y7").Value = 1
$ws.Range("Y8").Value = 18
$ws.Range("Y10").Value = 0
$ws.Range("Y11").Value = 0
$ws.Range("Y12").Value = 0
$ws.Range("Y13").Value = 0
$ws.Range("Y14").Value = 0
$ws.Range("Y15").Value = 0
$ws.Range("Y16").Value = 0
$ws.Range("Y17").Value = 0
$ws.Range("Y19").Value = 0
$ws.Range("Y21").Value = 0
$ws.Range("Y22").Value = 0
$ws.Range("Y23").Value = 0
$ws.Range("Y24").Value = 3
$ws.Range("Y27").Value = 0
$ws.Range("Y30").Value = 0
$ws.Range("Y31").Value = 0
$ws.Range("Y33").Value = 0
$ws.Range("Y34").Value = 18
$ws.Range("Y36").Value = 0
$ws.Range("Y37").Value = 0
$ws.Range("Y38").Value = 0
$ws.Range("Y39").Value = 0
$ws.Range("Y40").Value = 0
$ws.Range("Y41").Value = 0
$ws.Range("Y42").Value = 0
$ws.Range("Y43").Value = 0
$ws.Range("Y44").Value = 0
$ws.Range("Y45").Value = 0
$ws.Range("Y46").Value = 0
$ws.Range("Y47").Value = 0
$ws.Range("Y48").Value = 0
$ws.Range("Y49").Value = 0
$ws.Range("Y51").Value = 1
$ws.Range("Y52").Value = 0
$ws.Range("Y53").Value = 0
$ws.Range("Y54").Value = 0
$ws.Range("Y55").Value = 0
$ws.Range("Y56").Value = 0
